$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cached enum render values to their uppercased form
$ws.Range("A2").Value = "BUILDYOURDREAM"
$ws.Range("B2").Value = "ANTIQUEWHITE"

# Widen columns A and B to match the new (longer, uppercased) cached values
$ws.Columns.Item(1).ColumnWidth = 19.7142857142857
$ws.Columns.Item(2).ColumnWidth = 15.7142857142857

